$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 973.55554
$ws.Range("J17").Value = 973.55554
$ws.Range("L17").Value = 2920.66662
$ws.Range("N17").Value = -3256.66662
$ws.Range("H69").Value = 4160
$ws.Range("J69").Value = 4160
$ws.Range("L69").Value = 12480
$ws.Range("N69").Value = -14228
$ws.Range("H72").Value = 4160
$ws.Range("J72").Value = 4160
$ws.Range("L72").Value = 37440
$ws.Range("N72").Value = -46176
$ws.Range("H76").Value = 2424.5945
$ws.Range("I76").Value = 2415
$ws.Range("J76").Value = 2442.3076
$ws.Range("K76").Value = 2415
$ws.Range("L76").Value = 2442.3076
$ws.Range("M76").Value = -2100
$ws.Range("N76").Value = -3072.3076
$ws.Range("H79").Value = 2424.5945
$ws.Range("I79").Value = 2415
$ws.Range("J79").Value = 2442.3076
$ws.Range("K79").Value = 2415
$ws.Range("L79").Value = 2442.3076
$ws.Range("M79").Value = -1323
$ws.Range("N79").Value = -4626.3076
$ws.Range("H93").Value = 37933.332
$ws.Range("J93").Value = 37933.332
$ws.Range("L93").Value = 37933.332
$ws.Range("N93").Value = -42925.332
$ws.Range("H121").Value = 1550
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 2300
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 6900
$ws.Range("M121").Value = -653
$ws.Range("N121").Value = -10394
$ws.Range("H131").Value = 2467.2222
$ws.Range("I131").Value = 2243.5715
$ws.Range("J131").Value = 3250
$ws.Range("K131").Value = 6730.7145
$ws.Range("L131").Value = 9750
$ws.Range("M131").Value = -1690.7145
$ws.Range("N131").Value = -19830
$ws.Range("H141").Value = 1883.2151
$ws.Range("I141").Value = 1074.2048
$ws.Range("K141").Value = 3222.6144
$ws.Range("M141").Value = 1957.3856
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1275.5483
$ws.Range("I61").Value = 1108.9474
$ws.Range("J61").Value = 1539.3334
$ws.Range("K61").Value = 1108.9474
$ws.Range("L61").Value = 1539.3334
$ws.Range("M61").Value = -896.9474
$ws.Range("N61").Value = -1963.3334
$ws.Range("H74").Value = 18519598
$ws.Range("I74").Value = 26316830
$ws.Range("J74").Value = 1176.75
$ws.Range("K74").Value = 26316830
$ws.Range("L74").Value = 1176.75
$ws.Range("M74").Value = -26315956
$ws.Range("N74").Value = -2924.75
$ws.Range("H77").Value = 18519598
$ws.Range("I77").Value = 26316830
$ws.Range("J77").Value = 1176.75
$ws.Range("K77").Value = 131584150
$ws.Range("L77").Value = 5883.75
$ws.Range("M77").Value = -131579782
$ws.Range("N77").Value = -14619.75
$ws.Range("H102").Value = 1303.0435
$ws.Range("I102").Value = 1303.0435
$ws.Range("K102").Value = 1303.0435
$ws.Range("M102").Value = 318.9565
$ws.Range("H132").Value = 2584.0833
$ws.Range("I132").Value = 1464.7858
$ws.Range("J132").Value = 4151.1
$ws.Range("K132").Value = 4394.357400000001
$ws.Range("L132").Value = 12453.3
$ws.Range("M132").Value = -1864.357400000001
$ws.Range("N132").Value = -17513.3
$ws.Range("H134").Value = 23200
$ws.Range("I134").Value = 23200
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 23200
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -18130
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 1275.5483
$ws.Range("I136").Value = 1108.9474
$ws.Range("J136").Value = 1539.3334
$ws.Range("K136").Value = 3326.8422
$ws.Range("L136").Value = 4618.0002
$ws.Range("M136").Value = -776.8422
$ws.Range("N136").Value = -9718.0002
$ws.Range("H137").Value = 27800
$ws.Range("J137").Value = 27800
$ws.Range("L137").Value = 27800
$ws.Range("N137").Value = -38000
$ws.Range("H139").Value = 38625
$ws.Range("J139").Value = 38625
$ws.Range("L139").Value = 38625
$ws.Range("N139").Value = -48905
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1087.7727
$ws.Range("I94").Value = 564.4375
$ws.Range("J94").Value = 2483.3333
$ws.Range("K94").Value = 564.4375
$ws.Range("L94").Value = 2483.3333
$ws.Range("M94").Value = -113.4375
$ws.Range("N94").Value = -3385.3333
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 775.8182
$ws.Range("I131").Value = 372.15384
$ws.Range("J131").Value = 945.0968
$ws.Range("K131").Value = 1116.46152
$ws.Range("L131").Value = 2835.2904
$ws.Range("M131").Value = 3923.53848
$ws.Range("N131").Value = -12915.2904
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1879.7142
$ws.Range("I122").Value = 1683.8572
$ws.Range("J122").Value = 2271.4285
$ws.Range("K122").Value = 5051.571599999999
$ws.Range("L122").Value = 6814.2855
$ws.Range("M122").Value = -2601.571599999999
$ws.Range("N122").Value = -11714.2855
$ws.Range("H132").Value = 4023.4092
$ws.Range("I132").Value = 4123.2646
$ws.Range("J132").Value = 3683.9
$ws.Range("K132").Value = 12369.7938
$ws.Range("L132").Value = 11051.7
$ws.Range("M132").Value = -9839.793800000001
$ws.Range("N132").Value = -16111.7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2239.9092
$ws.Range("I61").Value = 1689.9231
$ws.Range("J61").Value = 3034.3333
$ws.Range("K61").Value = 1689.9231
$ws.Range("L61").Value = 3034.3333
$ws.Range("M61").Value = -1487.9231
$ws.Range("N61").Value = -3438.3333
$ws.Range("H100").Value = 1735.9286
$ws.Range("I100").Value = 1237.875
$ws.Range("K100").Value = 1237.875
$ws.Range("M100").Value = -696.875
$ws.Range("H113").Value = 2239.9092
$ws.Range("I113").Value = 1689.9231
$ws.Range("J113").Value = 3034.3333
$ws.Range("K113").Value = 1689.9231
$ws.Range("L113").Value = 3034.3333
$ws.Range("M113").Value = 480.0769
$ws.Range("N113").Value = -7374.3333
$ws.Range("H132").Value = 1379.2982
$ws.Range("I132").Value = 984.12244
$ws.Range("J132").Value = 3799.75
$ws.Range("K132").Value = 2952.36732
$ws.Range("L132").Value = 11399.25
$ws.Range("M132").Value = -422.3673199999998
$ws.Range("N132").Value = -16459.25
$ws.Range("H136").Value = 2229.9424
$ws.Range("I136").Value = 1789.9032
$ws.Range("J136").Value = 2879.524
$ws.Range("K136").Value = 5369.7096
$ws.Range("L136").Value = 8638.572
$ws.Range("M136").Value = -2819.7096
$ws.Range("N136").Value = -13738.572
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 73483.836
$ws.Range("J80").Value = 73483.836
$ws.Range("L80").Value = 73483.836
$ws.Range("N80").Value = -75479.836
$ws.Range("H83").Value = 73483.836
$ws.Range("J83").Value = 73483.836
$ws.Range("L83").Value = 220451.508
$ws.Range("N83").Value = -230435.508
$ws.Range("H96").Value = 2092.9473
$ws.Range("I96").Value = 1457.5714
$ws.Range("J96").Value = 2463.5833
$ws.Range("K96").Value = 1457.5714
$ws.Range("L96").Value = 2463.5833
$ws.Range("M96").Value = -84.57140000000004
$ws.Range("N96").Value = -5209.5833
$ws.Range("H136").Value = 4116.0586
$ws.Range("I136").Value = 1213.3
$ws.Range("J136").Value = 8262.857
$ws.Range("K136").Value = 3639.9
$ws.Range("L136").Value = 24788.571
$ws.Range("M136").Value = -1089.9
$ws.Range("N136").Value = -29888.571
